$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-07-08 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-09 Sunday", 2) | Out-Null

# Update each table cell value (row-major order matches the diff)
$t = $d.Tables.Item(1)
$values = @(
    "9+13=22",
    "41+20=61",
    "45+9=54",
    "14+79=93",
    "30+23=53",
    "38-37=1",
    "88-16=72",
    "32-14=18",
    "23+24=47",
    "56+36=92",
    "84+14=98",
    "64+11=75",
    "26+72=98",
    "69-15=54",
    "80-78=2",
    "45+49=94",
    "8+58=66",
    "0+11=11",
    "28-0=28",
    "10+11=21",
    "95-86=9",
    "53-7=46",
    "47-28=19",
    "98-10=88",
    "63-29=34",
    "46+6=52",
    "52-4=48",
    "23+49=72",
    "79-10=69",
    "97-1=96",
    "2+3=5",
    "93-41=52",
    "17+6=23",
    "41+46=87",
    "46-38=8",
    "45+34=79",
    "85-41=44",
    "26+61=87",
    "55-25=30",
    "91-75=16",
    "12+70=82",
    "14+81=95",
    "94-33=61",
    "51-44=7",
    "44+39=83",
    "68+16=84",
    "57-4=53",
    "11+49=60",
    "57-36=21",
    "24+47=71",
    "88+6=94",
    "40-16=24",
    "5+60=65",
    "3+14=17",
    "44+28=72",
    "13+61=74",
    "8+2=10",
    "57+0=57",
    "71-39=32",
    "64-62=2",
    "56+0=56",
    "56+29=85",
    "21+37=58",
    "37-25=12",
    "95-32=63",
    "70-66=4",
    "52-13=39",
    "45+39=84",
    "82-62=20",
    "62-23=39",
    "21+31=52",
    "91-48=43",
    "29+56=85",
    "87-6=81",
    "45-8=37",
    "21+48=69",
    "61+11=72",
    "47+46=93",
    "77-72=5",
    "9+1=10",
    "66+1=67",
    "2+10=12",
    "7+62=69",
    "25+41=66",
    "49+40=89",
    "16+78=94",
    "85-16=69",
    "56+20=76",
    "60-55=5",
    "13+16=29",
    "65-37=28",
    "16+18=34",
    "83-22=61",
    "95-65=30",
    "48-13=35",
    "93-88=5",
    "80-10=70",
    "24+58=82",
    "64-8=56",
    "90-40=50"
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Updated" $idx "cells"
